$d = $word.ActiveDocument

function Replace-ExactText($findText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $findText"
        return
    }
    $rng.Text = $newText
}

# Title
Replace-ExactText "Quantum Entanglement: Unveiling the Mysteries of Interconnectedness" "Mathematics: The Universal Language of Science and Logic"

# Author name
Replace-ExactText "Samuel Davies" "Jasper Lancaster"

# Author email -> single generic run
Replace-ExactText "samuel.davies@quantumstudies.org" "yourvalidname"

# Intro paragraph (3 sentences)
Replace-ExactText "Amidst the perplexing wonders of the quantum realm exists a profound phenomenon challenging our understanding of reality - quantum entanglement" "Mathematics, a subject that has fascinated and challenged minds for centuries, unveils the mysteries of the universe through the intricate tapestry of numbers, equations, and geometric patterns"

Replace-ExactText " This enigmatic connection between particles, regardless of their distance, has captured the imagination of scientists, philosophers, and artists alike" " Mathematics provides a universal language that transcends cultures, enabling us to understand the cosmos, unravel the enigmas of nature, and harness the power of logic to solve complex problems"

Replace-ExactText " In this essay, we embark on a journey to unravel the mysteries of quantum entanglement, exploring its implications for our comprehension of the universe and delving into the potential applications that may revolutionize various fields" " The field of mathematics is a symphony of abstract concepts, where symbols dance in harmony, revealing the underlying order and beauty of our world"

# EPR paragraph
Replace-ExactText "In 1935, Albert Einstein, Boris Podolsky, and Nathan Rosen introduced the concept of quantum entanglement through their famous thought experiment known as the EPR paradox" "Through the exploration of mathematical concepts, we unlock the secrets of nature's blueprint"

Replace-ExactText " Their proposal demonstrated that two particles, once entangled, remain interconnected regardless of the distance separating them" " The Fibonacci sequence, found in the spirals of seashells or the patterns of plant growth, exemplifies the intricate relationship between numbers and biological structures"

Replace-ExactText " This relationship transcends the constraints of space and time, allowing one particle to instantaneously influence the other, even across vast cosmological distances" " The elegance of geometric shapes, such as fractals, reflects the self-similarity found in everything from snowflakes to coastlines. These patterns underscore the profound interconnectedness of all things and provide a glimpse into the underlying mathematical principles that govern our universe"

# Moreover paragraph
Replace-ExactText "Moreover, quantum entanglement defies classical intuition" "Mathematics isn't merely a collection of abstract theories; it's a powerful tool with practical applications in every field imaginable"

Replace-ExactText " When entangled particles are measured, their properties, such as spin or polarization, are correlated in a way that cannot be explained by classical physics" " It empowers engineers to design structures that withstand earthquakes, enables us to predict weather patterns, and makes it possible to develop new medical treatments and technologies"

Replace-ExactText " This non-locality, as it is known, challenges our conventional notions of causality and raises fundamental questions about the nature of reality itself" " From the economy to finance, from computer science to data analysis, and even in music and art, the profound influence of mathematics is undeniable. It's a subject that touches every aspect of our lives, shaping our understanding of the world and guiding us towards a future filled with infinite possibilities"

# Summary heading paragraph
Replace-ExactText "Quantum entanglement, an awe-inspiring phenomenon, offers a glimpse into the uncharted territory of the quantum world" "Mathematics serves as an essential tool for understanding the intricacies of our universe, providing a lens through which we unlock the mysteries of science and logic"

Replace-ExactText " Its non-local nature challenges our fundamental understanding of reality, while its potential applications hold promise for transformative technologies" " Its abstract concepts find practical applications in diverse fields, empowering engineers, scientists, musicians, artists, and countless other professionals to innovate and drive progress"

Replace-ExactText " From quantum computing to secure communication, entanglement-based technologies may revolutionize numerous fields" " Mathematics unveils the interconnectedness of all things, from the Fibonacci sequence found in nature to the intricate patterns in art and music"

Replace-ExactText " Though much remains unknown, continued exploration of quantum entanglement promises to deepen our comprehension of the universe and expand the boundaries of human knowledge" " It's a subject that permeates our existence, shaping our understanding of the cosmos and enabling us to chart a course toward a future of endless possibilities"

# Add a new trailing empty paragraph at the very end of the document body
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "All replacements applied."
